{"js": "// Bele\u00edrtam a docx-be mint Ujcsi\n//\n// The document body is a single empty paragraph. Replace it with a\n// paragraph containing \"Ujcsi teszt\", split into two runs around the\n// proofErr spell-check markers (\"Ujcsi\" is flagged by the spell checker,\n// \"spellStart\"/\"spellEnd\" bracket it) exactly as Word records when you\n// type a word it doesn't recognize.\n//\n// The plain text APIs (body.insertText / paragraph.insertText) only ever\n// produce a single run, so we insert raw OOXML instead - insertOoxml\n// expects the WordprocessingML wrapped in the \"Flat OPC\" <pkg:package>\n// envelope, and with InsertLocation.replace it swaps out the whole body\n// (the lone empty paragraph) for the new paragraph.\n\nconst body = context.document.body;\n\nconst wordXml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n      '<w:p>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Ujcsi</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> teszt</w:t></w:r>' +\n      '</w:p>' +\n    '</w:body>' +\n  '</w:document>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' + wordXml + '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Bele\u00edrtam a docx-be mint Ujcsi\n#\n# The document has a single empty paragraph. Replace its content with two\n# runs - \"Ujcsi\" and \" teszt\" - separated by the proofErr spell-check\n# markers that Word inserts around a word it flags as a possible spelling\n# mistake while the user types (\"Ujcsi\" is not a dictionary word).\n#\n# Range.InsertXML expects a WordprocessingML package wrapped in the\n# \"Flat OPC\" <pkg:package> envelope; InsertXML replaces the contents of the\n# range it is called on, so calling it on the whole-document Content range\n# swaps out the lone empty paragraph for the new one.\n\n$d = $word.ActiveDocument\n\n$wordXml = '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n    '<w:p>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>Ujcsi</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> teszt</w:t></w:r>' +\n    '</w:p>' +\n  '</w:body>' +\n'</w:document>'\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' + $wordXml + '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$d.Content.InsertXML($flatOpc)\n"}
